$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$tbl.Cell(1, 1).Range.Text = "375×8=3000"
$tbl.Cell(1, 2).Range.Text = "252×8=2016"
$tbl.Cell(1, 3).Range.Text = "644×6=3864"
$tbl.Cell(1, 4).Range.Text = "845×8=6760"
$tbl.Cell(1, 5).Range.Text = "303×9=2727"

$tbl.Cell(5, 1).Range.Text = "813×7=5691"
$tbl.Cell(5, 2).Range.Text = "233×6=1398"
$tbl.Cell(5, 3).Range.Text = "853×6=5118"
$tbl.Cell(5, 4).Range.Text = "764×5=3820"
$tbl.Cell(5, 5).Range.Text = "872×4=3488"

$tbl.Cell(10, 1).Range.Text = "854×7=5978"
$tbl.Cell(10, 2).Range.Text = "920×3=2760"
$tbl.Cell(10, 3).Range.Text = "420×2=840"
$tbl.Cell(10, 4).Range.Text = "620×4=2480"
$tbl.Cell(10, 5).Range.Text = "915×8=7320"

$tbl.Cell(15, 1).Range.Text = "530×8=4240"
$tbl.Cell(15, 2).Range.Text = "686×6=4116"
$tbl.Cell(15, 3).Range.Text = "877×4=3508"
$tbl.Cell(15, 4).Range.Text = "341×7=2387"
$tbl.Cell(15, 5).Range.Text = "498×7=3486"

$tbl.Cell(20, 1).Range.Text = "169×3=507"
$tbl.Cell(20, 2).Range.Text = "618×2=1236"
$tbl.Cell(20, 3).Range.Text = "841×4=3364"
$tbl.Cell(20, 4).Range.Text = "472×7=3304"
$tbl.Cell(20, 5).Range.Text = "293×6=1758"
